$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 38971
$ws.Range("D2").Value = 56343594
$ws.Range("C3").Value = 93382
$ws.Range("D3").Value = 136872943
$ws.Range("C4").Value = 31882
$ws.Range("D4").Value = 47211967
$ws.Range("C5").Value = 8973
$ws.Range("D5").Value = 13336887
$ws.Range("C6").Value = 2086
$ws.Range("D6").Value = 3100971
$ws.Range("C7").Value = 172
$ws.Range("D7").Value = 253093
$ws.Range("C12").Value = 42346
$ws.Range("D12").Value = 57404469
$ws.Range("C13").Value = 9938
$ws.Range("D13").Value = 14369023
$ws.Range("C14").Value = 26531
$ws.Range("D14").Value = 38889576
$ws.Range("C15").Value = 8486
$ws.Range("D15").Value = 12594478
$ws.Range("C16").Value = 2214
$ws.Range("D16").Value = 3289539
$ws.Range("C17").Value = 432
$ws.Range("D17").Value = 637123
$ws.Range("C20").Value = 10440
$ws.Range("D20").Value = 13808950
$ws.Range("C21").Value = 13734
$ws.Range("D21").Value = 19819964
$ws.Range("C22").Value = 32340
$ws.Range("D22").Value = 47448771
$ws.Range("C23").Value = 10449
$ws.Range("D23").Value = 15530516
$ws.Range("C24").Value = 2708
$ws.Range("D24").Value = 4026771
$ws.Range("C25").Value = 544
$ws.Range("D25").Value = 810092
$ws.Range("C27").Value = 11953
$ws.Range("D27").Value = 15953633
$ws.Range("C28").Value = 7895
$ws.Range("D28").Value = 11424732
$ws.Range("C29").Value = 23067
$ws.Range("D29").Value = 33861815
$ws.Range("C30").Value = 7969
$ws.Range("D30").Value = 11852392
$ws.Range("C31").Value = 2016
$ws.Range("D31").Value = 3008251
$ws.Range("C32").Value = 379
$ws.Range("D32").Value = 565915
$ws.Range("C34").Value = 8507
$ws.Range("D34").Value = 11234675
$ws.Range("C35").Value = 3374
$ws.Range("D35").Value = 4873000
$ws.Range("C36").Value = 8055
$ws.Range("D36").Value = 11764926
$ws.Range("C37").Value = 3239
$ws.Range("D37").Value = 4801461
$ws.Range("C41").Value = 2543
$ws.Range("D41").Value = 3435728
$ws.Range("C42").Value = 17792
$ws.Range("D42").Value = 25726925
$ws.Range("C43").Value = 52362
$ws.Range("D43").Value = 76743814
$ws.Range("C44").Value = 19360
$ws.Range("D44").Value = 28749757
$ws.Range("C45").Value = 5755
$ws.Range("D45").Value = 8566187
$ws.Range("C46").Value = 1264
$ws.Range("D46").Value = 1886545
$ws.Range("C50").Value = 17182
$ws.Range("D50").Value = 22824514
$ws.Range("C51").Value = 2161
$ws.Range("D51").Value = 3136728
$ws.Range("C52").Value = 7271
$ws.Range("D52").Value = 10685390
$ws.Range("C53").Value = 2441
$ws.Range("D53").Value = 3645572
$ws.Range("C54").Value = 772
$ws.Range("D54").Value = 1153305
$ws.Range("C56").Value = 21
$ws.Range("D56").Value = 31500
$ws.Range("C57").Value = 7373
$ws.Range("D57").Value = 10140856
$ws.Range("C58").Value = 1201
$ws.Range("D58").Value = 2068952
$ws.Range("C59").Value = 2915
$ws.Range("D59").Value = 4981483
$ws.Range("C60").Value = 1149
$ws.Range("D60").Value = 1974106
$ws.Range("C61").Value = 396
$ws.Range("D61").Value = 684883
$ws.Range("C62").Value = 133
$ws.Range("D62").Value = 238100
$ws.Range("C64").Value = 1742
$ws.Range("D64").Value = 2762906
$ws.Range("C65").Value = 15871
$ws.Range("D65").Value = 22920177
$ws.Range("C66").Value = 45873
$ws.Range("D66").Value = 67110547
$ws.Range("C67").Value = 16034
$ws.Range("D67").Value = 23823506
$ws.Range("C68").Value = 4668
$ws.Range("D68").Value = 6951788
$ws.Range("C69").Value = 971
$ws.Range("D69").Value = 1444668
$ws.Range("C70").Value = 82
$ws.Range("D70").Value = 120330
$ws.Range("C73").Value = 15446
$ws.Range("D73").Value = 20341197
$ws.Range("C74").Value = 55005
$ws.Range("D74").Value = 80041560
$ws.Range("C75").Value = 153437
$ws.Range("D75").Value = 226017556
$ws.Range("C76").Value = 66042
$ws.Range("D76").Value = 98402314
$ws.Range("C77").Value = 21188
$ws.Range("D77").Value = 31659884
$ws.Range("C78").Value = 5080
$ws.Range("D78").Value = 7588403
$ws.Range("C85").Value = 53932
$ws.Range("D85").Value = 73253536
$ws.Range("C86").Value = 4792
$ws.Range("D86").Value = 6941864
$ws.Range("C87").Value = 11929
$ws.Range("D87").Value = 17521062
$ws.Range("C88").Value = 3977
$ws.Range("D88").Value = 5925458
$ws.Range("C93").Value = 5580
$ws.Range("D93").Value = 7497056
$ws.Range("C94").Value = 1665
$ws.Range("D94").Value = 2399699
$ws.Range("C95").Value = 5362
$ws.Range("D95").Value = 7898519
$ws.Range("C96").Value = 1990
$ws.Range("D96").Value = 2962926
$ws.Range("C97").Value = 709
$ws.Range("D97").Value = 1062460
$ws.Range("C98").Value = 198
$ws.Range("D98").Value = 297613
$ws.Range("C101").Value = 3707
$ws.Range("D101").Value = 4911213
$ws.Range("C102").Value = 739
$ws.Range("D102").Value = 1259775
$ws.Range("C103").Value = 460
$ws.Range("D103").Value = 816027
$ws.Range("C104").Value = 171
$ws.Range("D104").Value = 303180
$ws.Range("C105").Value = 55
$ws.Range("D105").Value = 96000
$ws.Range("C106").Value = 30
$ws.Range("D106").Value = 58500
$ws.Range("C107").Value = 11123
$ws.Range("D107").Value = 16132574
$ws.Range("C108").Value = 29846
$ws.Range("D108").Value = 43829317
$ws.Range("C109").Value = 10003
$ws.Range("D109").Value = 14872038
$ws.Range("C110").Value = 2761
$ws.Range("D110").Value = 4116580
$ws.Range("C112").Value = 57
$ws.Range("D112").Value = 85500
$ws.Range("C113").Value = 8
$ws.Range("D113").Value = 12000
$ws.Range("C114").Value = 10022
$ws.Range("D114").Value = 13230962
$ws.Range("C115").Value = 31335
$ws.Range("D115").Value = 45175351
$ws.Range("C116").Value = 67703
$ws.Range("D116").Value = 99058377
$ws.Range("C117").Value = 21797
$ws.Range("D117").Value = 32388076
$ws.Range("C118").Value = 6196
$ws.Range("D118").Value = 9229841
$ws.Range("C119").Value = 1163
$ws.Range("D119").Value = 1738100
$ws.Range("C120").Value = 87
$ws.Range("D120").Value = 126895
$ws.Range("C124").Value = 26407
$ws.Range("D124").Value = 35237741
$ws.Range("C125").Value = 37186
$ws.Range("D125").Value = 53655252
$ws.Range("C126").Value = 78811
$ws.Range("D126").Value = 115227889
$ws.Range("C127").Value = 24370
$ws.Range("D127").Value = 36170541
$ws.Range("C128").Value = 6547
$ws.Range("D128").Value = 9729123
$ws.Range("C129").Value = 1300
$ws.Range("D129").Value = 1933311
$ws.Range("C133").Value = 32562
$ws.Range("D133").Value = 43214787
$ws.Range("C134").Value = 13674
$ws.Range("D134").Value = 19793715
$ws.Range("C135").Value = 33081
$ws.Range("D135").Value = 48580447
$ws.Range("C136").Value = 11715
$ws.Range("D136").Value = 17405816
$ws.Range("C137").Value = 3036
$ws.Range("D137").Value = 4525241
$ws.Range("C141").Value = 11060
$ws.Range("D141").Value = 14739961
$ws.Range("C142").Value = 36303
$ws.Range("D142").Value = 52430055
$ws.Range("C143").Value = 83724
$ws.Range("D143").Value = 122653153
$ws.Range("C144").Value = 24951
$ws.Range("D144").Value = 37066958
$ws.Range("C145").Value = 6552
$ws.Range("D145").Value = 9776496
$ws.Range("C146").Value = 1496
$ws.Range("D146").Value = 2226230
$ws.Range("C147").Value = 87
$ws.Range("D147").Value = 130130
$ws.Range("C149").Value = 29959
$ws.Range("D149").Value = 40389868
